$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new entry was recorded for the "Others" group's September log: insert a
# fresh row above the existing first entry (row 30), which pushes every
# subsequent row (through the old row 67, the "Broadband" group header) down
# by one — growing the sheet's used range from A1:Y67 to A1:Y68.
$ws.Rows.Item(30).Insert()

# Populate the new row with the latest September_Details / September_Date
# values (columns R/S); all other columns on this row stay blank.
$ws.Range("R30").Value = "internet bal axisbank"
$ws.Range("S30").Value = "2024-09-05 16:05:55"
